$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 4041.375937365727
$ws.Range("D3").Value = 4041.375937365723
$ws.Range("D4").Value = 4041.375937365723
$ws.Range("D6").Value = 11546.78839247327
$ws.Range("D7").Value = 11546.78839247327
$ws.Range("D9").Value = 2062.858975181512
$ws.Range("D10").Value = 2062.85897518151
$ws.Range("D11").Value = 30
$ws.Range("D12").Value = 1702.85897518151
$ws.Range("D13").Value = 30
$ws.Range("D17").Value = 24000
$ws.Range("D19").Value = 17835.43650000022
$ws.Range("D20").Value = 17835.43650000022
$ws.Range("D21").Value = 360
$ws.Range("D24").Value = 186262.6104983193
$ws.Range("D25").Value = 186262.6104983191
$ws.Range("D28").Value = 3725.252209966534
$ws.Range("D29").Value = 3725.252209966534
$ws.Range("D30").Value = 186262.6104983191
$ws.Range("D35").Value = 20217.59999999952
$ws.Range("D36").Value = 20217.59999999952
$ws.Range("D38").Value = -452.7745246309626
$ws.Range("D39").Value = -452.7745246309621
$ws.Range("D41").Value = 452.7745246309621
$ws.Range("D42").Value = 9055.49049261916
$ws.Range("D43").Value = 9055.49049261916
